$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-17 01:55:48"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
